# The deck currently uses the "Integral" theme (applied to the one slide
# master / design in this file). The author switched the design back to the
# default "Office Theme" palette (Design tab -> Office Theme colors).
#
# PowerPoint exposes the twelve theme colour roles (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) through Slide.ThemeColorScheme.Colors(i).RGB,
# in that exact order - which is the live, editable surface for the design
# actually applied to the slides/slide master. We rewrite every slot from
# the old "Integral" values to the standard "Office Theme" values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function New-RGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index : role      : target "Office Theme" colour
$officeThemeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $tcs.Colors($i).RGB = New-RGB $officeThemeColors[$i - 1]
}
